$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing numeric / string values
$ws.Range("H4").Value = 6
$ws.Range("K6").Value = "headphones"
$ws.Range("K7").Value = 15
$ws.Range("H8").Value = 5
$ws.Range("K8").Value = 1
$ws.Range("H12").Value = 8

# Add new Test Result cells in row 26
$ws.Range("D26").Value = "V"
$ws.Range("I26").Value = "V"
